$d = $word.ActiveDocument

# The room-printout title paragraph ("Gebäude #g / Raum #r") shrinks from
# 36pt (72 half-points) down to 32pt (64 half-points), and its paragraph
# mark (previously 30pt / 60 half-points) grows to match at 32pt so the
# whole heading line is now a uniform 64 half-points.
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.Font.Size = 32
$r.Font.SizeBi = 32
